$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-17 17:03:39"

$wsZhCn.Range("H2").Value = "2016-08-17 17:03:34"
$wsZhCn.Range("K2").Value = "2016-08-17 17:03:56"

$wsDeDe.Range("K2").Value = "2016-08-17 17:04:09"
